$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its original text representation
# (values like "1.000" or "0.000007796" must not be reinterpreted as numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.166.39'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.861.66'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '0.7083'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '240.99'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.3092'
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("E9").Value = '  -3.45%  '
$ws.Range("D10").Value = '24.68'
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").Value = '0.08352'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '1.853.67'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '5.185'
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("D14").Value = '0.7084'
$ws.Range("E14").Value = '  -2.87%  '
$ws.Range("D15").Value = '91.09'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '29.203.80'
$ws.Range("E16").Value = '  -0.87%  '
$ws.Range("D17").Value = '5.912'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '242.65'
$ws.Range("E18").Value = '  -2.15%  '
$ws.Range("D19").Value = '0.000007796'
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").Value = '2.112.21'
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D23").Value = '7.873'
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("D24").Value = '0.9996'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '0.1583'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").Value = '163.25'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '8.937'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").Value = '18.45'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '1.331'
$ws.Range("E29").Value = '  -2.26%  '
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").Value = '4.403'
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = '4.226'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").Value = '0.05141'
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("D34").Value = '0.7963'
$ws.Range("E34").Value = '  +9.45%  '
$ws.Range("D35").Value = '1.915'
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  -2.94%  '
$ws.Range("D37").Value = '2.679'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").Value = '2.692'
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").Value = '1.168.27'
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("D41").Value = '6.212'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").Value = '0.8892'
$ws.Range("E42").Value = '  -2.58%  '
$ws.Range("D43").Value = '72.82'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").Value = '2.012.83'
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").Value = '0.5200'
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("D48").Value = '1.772'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = '9.333'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").Value = '0.4272'
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("B51").Value = 'Frax'
$ws.Range("C51").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").Value = '0.9982'
$ws.Range("E51").Value = '  -0.44%  '
